# The Adyen test-data fixture had its "Gross Currency" / "Net Currency"
# column values changed from "EUR" to "USD" throughout the sheet
# (everything else in the diff is incidental re-save noise from the
# authoring tool, not a deliberate content change).
#
# Doing this as a sheet-wide Find & Replace (rather than per-cell writes)
# mirrors how the edit would actually have been made in Excel, and lets
# the shared-string table update once for every cell that held "EUR".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$replaced = $used.Replace("EUR", "USD", 1, 1, $true, $false, $false, $false)

# Best-effort: the workbook's tab-bar/horizontal-scrollbar split ratio
# also moved (615 -> 500) in the saved file.
try {
    $win = $wb.Windows.Item(1)
    $win.TabRatio = 500
} catch {
}
